# Effector inputs workbook update
# - Adds a new "Effector Type (0 = Fixed)" column (D) with header + sample values
# - Adds a new effector data row (row 6, with row 5 left blank) for a negative-position effector
# - Updates the active selection to B8 (matches the post-edit cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "Effector Type (0 = Fixed)"

# Fill in the "Effector Type" values for the two existing fully-populated rows
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0

# New effector entry on row 6 (row 5 intentionally left empty)
$ws.Range("A6").Value = -500
$ws.Range("B6").Value = -500
$ws.Range("C6").Value = 1000

# Size the new column to fit its (longer) header text (~20.33 "characters"
# once Excel's standard 5/6-character padding is added back on save)
$ws.Columns.Item(4).ColumnWidth = 19.5

# Leave the selection where the user ended up after entering the data
$ws.Range("B8").Select()
